$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the PriceChange (X) / UpDown (Y) columns that were left blank ---
$ws.Range("X3").Value = 0.19000099999999875
$ws.Range("Y3").Value = "Up"

# --- Row 4: new day of data appended to the sheet ---

# Date column needs the same date/time display style already used by A2:A3,
# so copy that formatting across before writing the value.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4").Value = 42641.890694444446

$ws.Range("B4").Value = 13            # ScoreFinal
$ws.Range("C4").Value = "Buy"         # Verdict
$ws.Range("D4").Value = 22            # totalSentiment
$ws.Range("E4").Value = 7164          # wordCount
$ws.Range("F4").Value = 974           # sentenceCount
$ws.Range("G4").Value = 64            # posWordPercentage
$ws.Range("H4").Value = 33            # negWordPercentage
$ws.Range("I4").Value = 99            # posPhrasePercentage
$ws.Range("J4").Value = 0             # negPhrasePercentage
$ws.Range("K4").Value = 9234          # ElapsedMs
$ws.Range("L4").Value = 178           # posWordCount
$ws.Range("M4").Value = 94            # negWordCount
$ws.Range("N4").Value = 3             # positivePhraseCount
$ws.Range("O4").Value = 0             # negativePhraseCount
$ws.Range("P4").Value = "Noun"        # Method

$ws.Range("Q4").Value = 65.63785237683328  # RSI
$ws.Range("R4").Value = 0.48               # PEG

# 200Moving% / 50Moving% share the percentage style used by S2:T3.
$ws.Range("S3").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S4").Value = 0.0685

$ws.Range("T3").Copy()
$ws.Range("T4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("T4").Value = -0.0619

$ws.Range("U4").Value = 2.27          # PriceBook
$ws.Range("V4").Value = "N/A"         # Dividend
$ws.Range("W4").Value = 0             # Bollinger
